{"js": "// The edit swaps the trailing reference numbers between two list items:\n//   \"...mahalle-sokak isimleri cekilecek (10)\"  ->  \"...cekilecek (9)\"\n//   \"...ilce-mahalle-sokak classlari (9)\"        ->  \"...classlari (10)\"\n// i.e. the \"10\" and \"9\" inside the parentheses are swapped between the two\n// paragraphs (the surrounding \" (\" and \")\" text stays put).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two paragraphs involved in the swap using stable surrounding text.\nlet paraWith10 = null; // currently ends with \" (10)\"\nlet paraWith9 = null;  // currently ends with \" (9)\"\n\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t.indexOf(\"Selenium otomasyonu ile mahalle-sokak isimleri cekilecek\") !== -1) {\n    paraWith10 = p;\n  } else if (t.indexOf(\"ilce-mahalle-sokak\") !== -1 && t.indexOf(\"classlari\") !== -1) {\n    paraWith9 = p;\n  }\n}\n\nif (!paraWith10 || !paraWith9) {\n  throw new Error(\"Could not locate the target paragraphs for the number swap.\");\n}\n\n// Swap through a placeholder so the \"10\" and \"9\" never collide mid-edit.\nconst PLACEHOLDER = \"\\u0001\";\n\nlet hits = paraWith10.search(\"10\", { matchCase: true });\nhits.load(\"items/text\");\nawait context.sync();\nif (hits.items.length === 0) {\n  throw new Error('Could not find \"10\" in the expected paragraph.');\n}\nhits.items[0].insertText(PLACEHOLDER, \"Replace\");\nawait context.sync();\n\nhits = paraWith9.search(\"9\", { matchCase: true });\nhits.load(\"items/text\");\nawait context.sync();\nif (hits.items.length === 0) {\n  throw new Error('Could not find \"9\" in the expected paragraph.');\n}\nhits.items[0].insertText(\"10\", \"Replace\");\nawait context.sync();\n\nhits = paraWith10.search(PLACEHOLDER, { matchCase: true });\nhits.load(\"items/text\");\nawait context.sync();\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the placeholder to finish the swap.\");\n}\nhits.items[0].insertText(\"9\", \"Replace\");\nawait context.sync();\n", "ps1": "# The edit swaps the trailing reference numbers between two list items:\n#   \"...mahalle-sokak isimleri cekilecek (10)\"  ->  \"...cekilecek (9)\"\n#   \"...ilce-mahalle-sokak classlari (9)\"        ->  \"...classlari (10)\"\n# i.e. the \"10\" and \"9\" inside the parentheses are swapped between the two\n# paragraphs (the surrounding \" (\" and \")\" text stays put).\n\n$d = $word.ActiveDocument\n\n$paraWith10 = $null  # paragraph currently ending in \" (10)\"\n$paraWith9 = $null   # paragraph currently ending in \" (9)\"\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Selenium otomasyonu ile mahalle-sokak isimleri cekilecek*\") {\n        $paraWith10 = $p\n    } elseif ($t -like \"*ilce-mahalle-sokak*classlari*\") {\n        $paraWith9 = $p\n    }\n}\n\nif ($paraWith10 -eq $null -or $paraWith9 -eq $null) {\n    throw \"Could not locate the target paragraphs for the number swap.\"\n}\n\n# Replace \"10\" -> \"9\", scoped to that paragraph's own Range only.\n$range1 = $paraWith10.Range\n$found1 = $range1.Find.Execute(\"10\", $false, $false, $false, $false, $false, $true, 1, $false, \"9\", 2)\nif (-not $found1) {\n    throw 'Could not find \"10\" in the expected paragraph.'\n}\n\n# Replace \"9\" -> \"10\", scoped to that paragraph's own Range only (unaffected\n# by the edit above since the ranges are independent).\n$range2 = $paraWith9.Range\n$found2 = $range2.Find.Execute(\"9\", $false, $false, $false, $false, $false, $true, 1, $false, \"10\", 2)\nif (-not $found2) {\n    throw 'Could not find \"9\" in the expected paragraph.'\n}\n"}
